# Test Case Added (Given When Then)

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1: header text tweak (TestCase -> TestScenarios per the commit) ---
$ws1.Range("B2").Value = "TestScenarios"
$ws1.Range("B13").Value = "TestScenarios"

# --- Sheet1: selection/view changes ---
$ws1.Range("B5").Select()

# --- Rename Sheet2 -> SwipeIn and build its content ---
$ws2.Name = "SwipeIn"

$ws2.Columns.Item(1).ColumnWidth = 10.875
$ws2.Columns.Item(2).ColumnWidth = 36.625
$ws2.Columns.Item(3).ColumnWidth = 34.5
$ws2.Columns.Item(4).ColumnWidth = 32.875

$ws2.Range("C2:C5").WrapText = $true
$ws2.Range("D3").WrapText = $true

$ws2.Range("A1:D1").Value = "Hdr"
$ws2.Range("A1").Value = "TestCase No"
$ws2.Range("B1").Value = "Given"
$ws2.Range("C1").Value = "When"
$ws2.Range("D1").Value = "Then"
$ws2.Range("A1:D1").Font.Bold = $true
$ws2.Range("A1:D1").Interior.Color = 65535

$ws2.Range("A2").Value = 1
$ws2.Range("B2").Value = "MetroCard is SwippedIn"
$ws2.Range("C2").Value = "When CardId has allowed characters `nAnd Balance is 5.5`nAnd StationId is A1`nAnd Day is Monday"
$ws2.Range("D2").Value = "User Able to CheckIn without Exception"

$ws2.Range("A3").Value = 2
$ws2.Range("B3").Value = "MetroCard is SwippedIn"
$ws2.Range("C3").Value = "When CardId has allowed characters `nAnd Balance is 5.4`nAnd StationId is A1`nAnd Day is Monday"
$ws2.Range("D3").Value = "User not Able to CheckIn and Exception is thrown `"Out of Balance`""

$ws2.Range("A4").Value = 3
$ws2.Range("B4").Value = "MetroCard is SwippedIn"
$ws2.Range("C4").Value = "When CardId has allowed characters `nAnd Balance is 5.6`nAnd StationId is A1`nAnd Day is Monday"
$ws2.Range("D4").Value = "User Able to CheckIn without Exception"

$ws2.Range("B5").Value = "MetroCard is SwippedIn"
$ws2.Range("C5").Value = "When CardId has allowed characters  `nAnd Balance is 5.6`nAnd StationId is A1`nAnd Day is Monday"
$ws2.Range("D5").Value = "User Able to CheckIn without Exception"

$ws2.Range("B4").Select()
$ws2.Activate()
